$wb = $excel.ActiveWorkbook

# --- Sheet 1 "Account Information": remove the blank formatted row 2 (it only
#     carried leftover Hyperlink/Date cell formatting, no real data) so every
#     account record shifts up one row ---
$ws1 = $wb.Worksheets.Item("Account Information")
$ws1.Rows.Item(2).Delete()
[void]$ws1.Range("F9").Select()

# --- Sheet 2 "Timeslot Information": close the gap — "Class ID" moves from
#     C1 into B1 ---
$ws2 = $wb.Worksheets.Item("Timeslot Information")
$ws2.Range("B1").Value = $ws2.Range("C1").Value2
$ws2.Range("C1").ClearContents()
$ws2.Activate()
[void]$ws2.Range("E6").Select()

# --- Drop the now-unused "Hyperlink" cell style left over from the removed
#     row (sheet no longer contains any hyperlink-formatted cells) ---
$wb.Styles.Item("Hyperlink").Delete()
